# Atualizacao de bases das ligas: swap the data (columns B:AD) between pairs
# of rows on the "Brazil Serie C" sheet. Column A (the running index) is left
# untouched; only the match data that was attached to it changes places.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each pair lists the two worksheet rows whose B:AD contents must be
# exchanged with one another.
$pairs = @(
    @(29,30),
    @(54,55),
    @(60,61),
    @(62,63),
    @(98,99),
    @(104,105),
    @(107,108),
    @(124,125),
    @(128,129),
    @(130,131),
    @(164,165)
)

# Columns B through AD (id, Div, Date, HomeTeam, ... PL_AhUnder)
$firstCol = 2   # B
$lastCol  = 30  # AD

foreach ($pair in $pairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    # Read every value of both rows first so the swap is not clobbered
    # while we write the new contents back.
    $valsA = @()
    $valsB = @()
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $valsA += , ($ws.Cells.Item($rowA, $col).Value2)
        $valsB += , ($ws.Cells.Item($rowB, $col).Value2)
    }

    $idx = 0
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        # Skip cells whose value is identical in both rows (e.g. Div/Date)
        # so we do not needlessly dirty cells that do not actually change.
        if ($valsA[$idx] -ne $valsB[$idx]) {
            $ws.Cells.Item($rowA, $col).Value2 = $valsB[$idx]
            $ws.Cells.Item($rowB, $col).Value2 = $valsA[$idx]
        }
        $idx++
    }
}
